$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Current (before) layout ---
# A: Date               B: Hearing time         C: Case reference number
# D: Case name          E: Judge(s)             F: Member(s)
# G: Mode of hearing    H: Venue                I: Additional information
#
# --- Target (after) layout ---
# A: Date               B: Hearing time         C: Case reference number
# D: Case name          E: Panel                F: Mode of hearing
# G: Venue              H: Additional information

# 1) Combine Judge(s)/Member(s) columns (E/F) into a single "Panel" value in column E
$panel2 = $ws.Range("E2").Value2.ToString() + ", " + $ws.Range("F2").Value2.ToString()
$panel3 = $ws.Range("E3").Value2.ToString() + ", " + $ws.Range("F3").Value2.ToString()

$ws.Range("E2").Value = $panel2
$ws.Range("E3").Value = $panel3

# 2) Delete the now-redundant "Member(s)" column (F). This shifts G/H/I left to F/G/H.
$ws.Columns.Item(6).Delete() | Out-Null

# After this delete:
# A: Date  B: Hearing time  C: Case reference number  D: Case name  E: Panel
# F: Mode of hearing  G: Venue  H: Additional information

# 3) Update header row text
$ws.Range("A1").Value = "Date"
$ws.Range("B1").Value = "Hearing time"
$ws.Range("C1").Value = "Case reference number"
$ws.Range("D1").Value = "Case name"
$ws.Range("E1").Value = "Panel"
$ws.Range("F1").Value = "Mode of hearing"
$ws.Range("G1").Value = "Venue"
$ws.Range("H1").Value = "Additional information"

# 4) Row 3: date stays a real date value; re-use A2's existing short-date style (same
#    number format already applied to A2, numFmtId 14) for A3 rather than its old
#    dd/mm/yyyy custom style, so no new number-format entry is created
$ws.Range("A2").Copy() | Out-Null
$ws.Range("A3").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$ws.Range("A3").Value = 45642
$excel.CutCopyMode = 0

# 5) Row 2: date becomes literal text "12/16/2024"; case reference becomes text "1234.0"
#    (force text interpretation, then clear the number format back to General so the
#    cell keeps using the default style, matching the authored workbook)
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "12/16/2024"
$ws.Range("A2").ClearFormats() | Out-Null

$ws.Range("C2").NumberFormat = "@"
$ws.Range("C2").Value = "1234.0"
$ws.Range("C2").ClearFormats() | Out-Null

$ws.Range("C3").NumberFormat = "@"
$ws.Range("C3").Value = "1235.0"
$ws.Range("C3").ClearFormats() | Out-Null

# 6) Selection / cursor position matches target file
$ws.Range("H7").Select() | Out-Null
